$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nädal 7")

# Row 15 (task 9)
$ws.Range("B15").Value = 43904
$ws.Range("C15").Value = 0.89583333333333337
$ws.Range("D15").Value = 0.97916666666666663
$ws.Range("F15").Value = 120
$ws.Range("G15").Value = "Raamatu lugemine"
$ws.Range("H15").Value = "Clean Code"
$ws.Range("J15").Value = "x"

# Row 16 (task 10)
$ws.Range("A16").Value = 10
$ws.Range("B16").Value = 43905
$ws.Range("C16").Value = 0.9375
$ws.Range("D16").Value = 0.98958333333333337
$ws.Range("F16").Value = 75
$ws.Range("G16").Value = "Raamatu lugemine"
$ws.Range("H16").Value = "Clean Code"
$ws.Range("J16").Value = "x"

# Row 17 (task 11)
$ws.Range("A17").Value = 11
$ws.Range("B17").Value = 43906
$ws.Range("C17").Value = 0.67361111111111116
$ws.Range("D17").Value = 0.71666666666666667
$ws.Range("F17").Value = 62
$ws.Range("G17").Value = "Kodutöö 7"
$ws.Range("H17").Value = "p. 28 tehtud"
$ws.Range("J17").Value = "x"

# Update the total formula to include the new rows
$ws.Range("F19").Formula = "=SUM(F7:F17)"

# Update the active selection to match the saved state
$ws.Range("J17").Select()
